# Update NATMI Jag2-Notch2 LR-pair sheet with new TPM-derived values.
#
# Layout of the sheet (rows 2..17, 20 data rows total): 16 rows arranged as
# 4 "sending cluster" blocks (column A) of 4 rows each; within each block the
# "target cluster" (column D) cycles through the same 4 clusters in the same
# order:
#   Sending clusters (outer, rows 2-5,6-9,10-13,14-17): ECs, FAPs, MuSCs, Resolving-Mac
#   Target  clusters (inner, within each block):        ECs, FAPs, MuSCs, Resolving-Mac
#
# Column meaning / dependency:
#   G = Ligand average expression value   -> depends only on sending cluster (A)
#   H = Ligand total   expression value   -> depends only on sending cluster (A)
#   M = Receptor average expression value -> depends only on target   cluster (D)
#   N = Receptor total  expression value  -> depends only on target   cluster (D)
#   I = G / sum(G across the 4 rows sharing the same target cluster D)
#   J = H / sum(H across the 4 rows sharing the same target cluster D)
#   O = M / sum(M across the 4 rows sharing the same sending cluster A)
#   P = N / sum(N across the 4 rows sharing the same sending cluster A)
#   Q = G * M
#   R = H * N
#   S = Q / sum(Q across all 16 rows)
#   T = R / sum(R across all 16 rows)
#
# New underlying (TPM-updated) per-cluster values, in cluster order
# ECs, FAPs, MuSCs, Resolving-Mac:

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")

# New ligand (average, total) expression values, keyed by sending cluster.
$newG = @{
    "ECs"           = 20.81798233333333
    "FAPs"          = 0.7925996666666667
    "MuSCs"         = 1.536855
    "Resolving-Mac" = 2.108971
}
$newH = @{
    "ECs"           = 62.453947
    "FAPs"          = 2.377799
    "MuSCs"         = 4.610564999999999
    "Resolving-Mac" = 6.326912999999999
}

# New receptor (average, total) expression values, keyed by target cluster.
$newM = @{
    "ECs"           = 1.910418
    "FAPs"          = 31.995262
    "MuSCs"         = 37.858701
    "Resolving-Mac" = 33.83466466666667
}
$newN = @{
    "ECs"           = 5.731254
    "FAPs"          = 95.985786
    "MuSCs"         = 113.576103
    "Resolving-Mac" = 101.503994
}

# Build row -> (sendingCluster, targetCluster) map matching the sheet layout
# (rows 2..17, 4 sending-cluster blocks of 4 target rows each).
$rowInfo = @{}
$r = 2
foreach ($send in $clusters) {
    foreach ($targ in $clusters) {
        $rowInfo[$r] = @{ Send = $send; Targ = $targ }
        $r = $r + 1
    }
}

$rows = 2..17

# --- Write G, H, M, N first ---
foreach ($row in $rows) {
    $info = $rowInfo[$row]
    $ws.Cells.Item($row, 7).Value2  = $newG[$info.Send]   # G
    $ws.Cells.Item($row, 8).Value2  = $newH[$info.Send]   # H
    $ws.Cells.Item($row, 13).Value2 = $newM[$info.Targ]   # M
    $ws.Cells.Item($row, 14).Value2 = $newN[$info.Targ]   # N
}

# --- Compute I, J (normalize G, H within rows sharing the same target cluster) ---
foreach ($targ in $clusters) {
    $targRows = @()
    foreach ($row in $rows) {
        if ($rowInfo[$row].Targ -eq $targ) { $targRows += $row }
    }

    $sumG = 0.0
    $sumH = 0.0
    foreach ($row in $targRows) {
        $sumG = $sumG + $newG[$rowInfo[$row].Send]
        $sumH = $sumH + $newH[$rowInfo[$row].Send]
    }

    foreach ($row in $targRows) {
        $g = $newG[$rowInfo[$row].Send]
        $h = $newH[$rowInfo[$row].Send]
        $ws.Cells.Item($row, 9).Value2  = $g / $sumG   # I
        $ws.Cells.Item($row, 10).Value2 = $h / $sumH   # J
    }
}

# --- Compute O, P (normalize M, N within rows sharing the same sending cluster) ---
foreach ($send in $clusters) {
    $sendRows = @()
    foreach ($row in $rows) {
        if ($rowInfo[$row].Send -eq $send) { $sendRows += $row }
    }

    $sumM = 0.0
    $sumN = 0.0
    foreach ($row in $sendRows) {
        $sumM = $sumM + $newM[$rowInfo[$row].Targ]
        $sumN = $sumN + $newN[$rowInfo[$row].Targ]
    }

    foreach ($row in $sendRows) {
        $m = $newM[$rowInfo[$row].Targ]
        $n = $newN[$rowInfo[$row].Targ]
        $ws.Cells.Item($row, 15).Value2 = $m / $sumM   # O
        $ws.Cells.Item($row, 16).Value2 = $n / $sumN   # P
    }
}

# --- Compute Q, R (edge weights) = G*M, H*N ---
$sumQ = 0.0
$sumR = 0.0
foreach ($row in $rows) {
    $info = $rowInfo[$row]
    $g = $newG[$info.Send]
    $h = $newH[$info.Send]
    $m = $newM[$info.Targ]
    $n = $newN[$info.Targ]
    $q = $g * $m
    $rVal = $h * $n
    $ws.Cells.Item($row, 17).Value2 = $q      # Q
    $ws.Cells.Item($row, 18).Value2 = $rVal   # R
    $sumQ = $sumQ + $q
    $sumR = $sumR + $rVal
}

# --- Compute S, T (normalize Q, R across the whole table) ---
foreach ($row in $rows) {
    $q = $ws.Cells.Item($row, 17).Value2
    $rVal = $ws.Cells.Item($row, 18).Value2
    $ws.Cells.Item($row, 19).Value2 = $q / $sumQ      # S
    $ws.Cells.Item($row, 20).Value2 = $rVal / $sumR   # T
}
